$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.527.63"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.080.48"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.60"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.53"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0780"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.385.96"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.41"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.74"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.784"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.22"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.066.16"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.451.88"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.72"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.00"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.60"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.86"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("E29").Value = "  -5.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.18"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.57"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.51"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.63"
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.492.21"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0958"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.32"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.25"
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.33"
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.26"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.98"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.271.86"
$ws.Range("E51").Value = "  +0.74%  "
